$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.845.53'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = '1.831.34'
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.85'
$ws.Range("E5").Value = '  -0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  +0.88%  '

$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.24'
$ws.Range("E8").Value = '  -2.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.327'
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0684'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0987'

$ws.Range("D12").Value = '2.093.99'
$ws.Range("E12").Value = '  +0.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.32'
$ws.Range("E13").Value = '  +2.02%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.668'
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.804.81'
$ws.Range("E15").Value = '  -0.64%  '

$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").Value = '34.803.71'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.49'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").Value = '0.0₃0787'
$ws.Range("E19").Value = '  -0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.10'
$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.16'
$ws.Range("E21").Value = '  +2.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.69'
$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.74'
$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.76'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.124'
$ws.Range("E27").Value = '  +2.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.32'
$ws.Range("E28").Value = '  -0.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.51'
$ws.Range("E29").Value = '  -7.48%  '

$ws.Range("E30").Value = '  +0.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0551'
$ws.Range("E31").Value = '  -0.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.91'
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.93'
$ws.Range("E33").Value = '  -1.13%  '

$ws.Range("E34").Value = '  +3.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.23'
$ws.Range("E35").Value = '  +8.43%  '

$ws.Range("E36").Value = '  +10.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.701'
$ws.Range("E37").Value = '  +2.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '91.39'
$ws.Range("E38").Value = '  -1.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.05'
$ws.Range("E39").Value = '  +6.46%  '

$ws.Range("D40").Value = '1.342.59'
$ws.Range("E40").Value = '  +2.80%  '

$ws.Range("E41").Value = '  -0.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.46'
$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("E43").Value = '  -1.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.25'
$ws.Range("E44").Value = '  -3.90%  '

$ws.Range("E45").Value = '  -0.44%  '

$ws.Range("E46").Value = '  -0.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0520'

$ws.Range("D48").Value = '2.009.52'
$ws.Range("E48").Value = '  +0.83%  '

$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0673'
$ws.Range("E50").Value = '  +4.38%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.24'
$ws.Range("E51").Value = '  -0.97%  '

